# The SOR based on greedy algorithm
#
# Updates the best bid/ask price ladder on the "Microsofts" sheet and
# switches the active tab / selection from "Netflix" back to "Microsofts".

$wb = $excel.ActiveWorkbook

$wsMicrosofts = $wb.Worksheets.Item("Microsofts")
$wsNetflix    = $wb.Worksheets.Item("Netflix")

# --- Price/quote updates on the "Microsofts" sheet ---------------------
$wsMicrosofts.Range("C4").Value  = 183.97

$wsMicrosofts.Range("A6").Value  = 183.92

$wsMicrosofts.Range("A7").Value  = 183.81
$wsMicrosofts.Range("C7").Value  = 183.85

$wsMicrosofts.Range("A12").Value = 183.73

$wsMicrosofts.Range("A14").Value = 183.71
$wsMicrosofts.Range("C14").Value = 183.78

# --- Active sheet / selection -------------------------------------------
# Move the active tab and selection back to "Microsofts" (away from
# "Netflix"), with the cursor parked just past the printed table.
$wsMicrosofts.Activate()
$wsMicrosofts.Range("E17").Select()
